$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new data row values
$ws.Range("A2").Value = "MCH108"
$ws.Range("C2").Value = "BOOK: A NOTE ON `"PLANNING`"THE UNION OF SOUTH AFRICA BY ANONYMOUS ADMIRER AFTER A MONTHS VISIT"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# Apply the font used for the new row (Calibri 10, theme text color)
# Resetting to the "Normal" cell style first yields a style without an
# explicit alignment override (matches A2/C2/D2/E2/G2/H2 in the target).
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.ThemeColor = 1

$ws.Range("C2:E2").Style = "Normal"
$ws.Range("C2:E2").Font.Name = "Calibri"
$ws.Range("C2:E2").Font.ThemeColor = 1

$ws.Range("G2:H2").Style = "Normal"
$ws.Range("G2:H2").Font.Name = "Calibri"
$ws.Range("G2:H2").Font.ThemeColor = 1

# F2 (the "1 Box" cell) carries a slightly different style (applyAlignment flag set)
$ws.Range("F2").Font.Name = "Calibri"
$ws.Range("F2").Font.ThemeColor = 1
